$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "Use rod..." note (row 2, column J): thuner -> thunder
$ws.Range("J2").Value = "Use rod to go rainbow: everything is radiant + thunder in paragon"

# Add new character row 6: Invoker / Controller-striker / Shardmind
$ws.Range("B6").Value = "Invoker"
$ws.Range("C6").Value = "Controller/striker"
$ws.Range("D6").Value = "Shardmind"

# Move the active selection down to B7, the next empty row
$ws.Range("B7").Select()
